$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 ---
$ws.Range("A14").Value = 11101100
$ws.Range("B14").Formula = "=LEN(A14)"
$ws.Range("C14").Formula = "=BIN2HEX(A14)"
$ws.Range("D14").Value = 2
$ws.Range("F14").Value = "Min"

# --- Row 15 ---
$ws.Range("A15").Value = 1110011
$ws.Range("D15").Value = 1
$ws.Range("F15").Value = 800000
$ws.Range("G15").Formula = "=HEX2DEC(F15)"

# --- Row 16 ---
$ws.Range("A16").Value = 11101101
$ws.Range("D16").Value = 0
$ws.Range("F16").Value = "Max"

# Shared formula B15:B20 (=LEN(A15)), then remove the B17 cell (no data there)
$ws.Range("B15:B20").Formula = "=LEN(A15)"
$ws.Range("B17").ClearContents()

# Shared formula C15:C16 (=BIN2HEX(A15))
$ws.Range("C15:C16").Formula = "=BIN2HEX(A15)"

# --- Row 17 ---
$ws.Range("F17").Value = "7FFFFF"
$ws.Range("G17").Formula = "=HEX2DEC(F17)"

# --- Row 18 ---
$ws.Range("A18").Value = 11101100
$ws.Range("C18").Formula = "=BIN2HEX(A18)"
$ws.Range("D18").Value = 2

# --- Row 19 ---
$ws.Range("A19").Value = 1111010
$ws.Range("D19").Value = 1
$ws.Range("H19").Value = [double]"1.11011000111001E+23"
$ws.Range("H19").NumberFormat = "0.00E+00"
$ws.Range("I19").Formula = "=BIN2DEC(H19)"

# --- Row 20 ---
$ws.Range("A20").Value = 1010000
$ws.Range("D20").Value = 0

# Shared formula C19:C20 (=BIN2HEX(A19))
$ws.Range("C19:C20").Formula = "=BIN2HEX(A19)"

# --- Row 22 ---
$ws.Range("A22").Value = 11101100
$ws.Range("C22").Formula = "=BIN2HEX(A22)"

# --- Row 23 ---
$ws.Range("A23").Value = 1010000

# --- Row 24 ---
$ws.Range("A24").Value = 110010

# Shared formula B22:B24 (=LEN(A22))
$ws.Range("B22:B24").Formula = "=LEN(A22)"

# Shared formula C23:C24 (=BIN2HEX(A23))
$ws.Range("C23:C24").Formula = "=BIN2HEX(A23)"

# --- Column widths ---
$ws.Columns.Item(6).ColumnWidth = 6.14
$ws.Columns.Item(8).ColumnWidth = 27.75

# --- View ---
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("H19").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
